# feat: add 2022-Q4 data
#
# The workbook currently has two sheets: "总计" (summary) and "2022-Q3"
# (fund holdings for Q3). We are adding a new quarter's data (2022-Q4):
#   - A new "2022-Q3" sheet is split off (a snapshot of the old data),
#     keeping the old Q3 numbers untouched.
#   - The original "2022-Q3" sheet is renamed to "2022-Q4" and refreshed
#     with the new quarter's fund numbers.
#   - The "总计" summary sheet gets a new row for the 2022-Q3 total, and
#     its former Q3 row now reports the Q4 totals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Duplicate the existing "2022-Q3" sheet so its current data survives
#    under its own tab, then repurpose the original tab for "2022-Q4".
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy([System.Reflection.Missing]::Value, $q3)

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

$newQ3 = $wb.Worksheets.Item(3)
$newQ3.Name = "2022-Q3"

# ---------------------------------------------------------------------
# 2) Refresh the "2022-Q4" sheet (the old "2022-Q3" tab) with the new
#    quarter's fund holdings. Fund codes / percentages are kept as text
#    (leading-zero codes like "012868" must not become numbers), so
#    values are entered with a leading apostrophe to force text, then the
#    cell style is reset to Normal so no stray quote-prefix format sticks.
# ---------------------------------------------------------------------
$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "'012868"
$q4.Range("C2").Value = "易方达标普信息科技指数（QDII-LOF）人民币 C"
$q4.Range("D2").Value = "'5.09"
$q4.Range("E2").Value = "'91.36"
$q4.Range("F2").Value = "'1.66"
$q4.Range("G2").Value = "'0.0845"
$q4.Range("H2").Value = 10

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "'161128"
$q4.Range("C3").Value = "易方达标普信息科技指数（QDII-LOF）人民币"
$q4.Range("D3").Value = "'5.09"
$q4.Range("E3").Value = "'91.36"
$q4.Range("F3").Value = "'1.66"
$q4.Range("G3").Value = "'0.0845"
$q4.Range("H3").Value = 10

$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "'003721"
$q4.Range("C4").Value = "易方达标普信息科技指数（QDII-LOF）美元A"
$q4.Range("D4").Value = "'4.93"
$q4.Range("E4").Value = "'91.36"
$q4.Range("F4").Value = "'1.66"
$q4.Range("G4").Value = "'0.0818"
$q4.Range("H4").Value = 10

$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "'012869"
$q4.Range("C5").Value = "易方达标普信息科技指数（QDII-LOF）美元 C"
$q4.Range("D5").Value = "'0.16"
$q4.Range("E5").Value = "'91.36"
$q4.Range("F5").Value = "'1.66"
$q4.Range("G5").Value = "'0.0027"
$q4.Range("H5").Value = 10

$q4.Range("B2:B5").Style = "Normal"
$q4.Range("D2:G5").Style = "Normal"

# The header row (B1:H1) and the A-column index cells (A2:A5) pick up the
# same "centered / bordered" style already used by the 总计 sheet's header
# (cellXf index 2), instead of the older style the original Q3 sheet used.
$summaryHeader = $wb.Worksheets.Item("总计").Range("B1")
$summaryHeader.Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats
$q4.Range("A2:A5").PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------
# 3) Update the "总计" summary sheet: the existing Q3 row now reflects
#    Q4, and a new row is appended with the original Q3 totals.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Range("B2").Value = "2022-Q4"

$summary.Range("A2").Copy($summary.Range("A3"))
$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 4
$summary.Range("D3").Value = 0.25
